# Updated cryptos list on Tue Aug  1 21:36:21 UTC 2023 with GitHub Actions
# Refreshes Price (col D) and Volume(1h) (col E) for every coin row, and for
# a few rows the scraped coin order shifted so Coin (col B) / Link (col C)
# are updated too. NumberFormat is forced to Text ("@") before assigning any
# Price value that Excel could otherwise auto-parse as a number, so values
# like "1.000" or "6.620" keep their exact original text instead of being
# coerced into 1 / 6.62.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.239.31"
$ws.Range("E2").Value = "  +0.03%  "

$ws.Range("D3").Value = "1.850.35"
$ws.Range("E3").Value = "  -0.30%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.76"
$ws.Range("E5").Value = "  +1.93%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7031"
$ws.Range("E6").Value = "  +0.42%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3074"
$ws.Range("E9").Value = "  -0.17%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.68"
$ws.Range("E10").Value = "  -0.45%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07829"
$ws.Range("E11").Value = "  +0.29%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "93.04"
$ws.Range("E12").Value = "  +0.90%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.140"
$ws.Range("E13").Value = "  +0.95%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.848.17"
$ws.Range("E14").Value = "  -0.53%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6878"
$ws.Range("E15").Value = "  +0.21%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.620"
$ws.Range("E16").Value = "  +1.78%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008326"
$ws.Range("E17").Value = "  -0.78%  "

$ws.Range("D18").Value = "29.226.70"
$ws.Range("E18").Value = "  -0.01%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.64"
$ws.Range("E19").Value = "  -2.82%  "

$ws.Range("D20").Value = "2.090.98"
$ws.Range("E20").Value = "  -1.19%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.76"
$ws.Range("E21").Value = "  -0.48%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9994"
$ws.Range("E22").Value = "  -0.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.533"
$ws.Range("E23").Value = "  +0.29%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9996"
$ws.Range("E24").Value = "  -0.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1517"
$ws.Range("E25").Value = "  +0.90%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.32"
$ws.Range("E26").Value = "  -0.73%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.849"
$ws.Range("E27").Value = "  +0.13%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.32"
$ws.Range("E28").Value = "  -0.84%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.541"
$ws.Range("E29").Value = "  -1.16%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.234"
$ws.Range("E30").Value = "  +0.09%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.188"
$ws.Range("E31").Value = "  +0.32%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.204"
$ws.Range("E32").Value = "  +0.84%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05123"
$ws.Range("E33").Value = "  -1.38%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7935"
$ws.Range("E34").Value = "  +4.37%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.927"
$ws.Range("E35").Value = "  +4.64%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.149"
$ws.Range("E36").Value = "  -1.23%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.696"
$ws.Range("E37").Value = "  -0.64%  "

$ws.Range("D38").Value = "1.324.87"
$ws.Range("E38").Value = "  +8.84%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01874"
$ws.Range("E39").Value = "  +0.67%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.717"
$ws.Range("E40").Value = "  -0.21%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9546"
$ws.Range("E41").Value = "  +6.62%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.072"
$ws.Range("E42").Value = "  +9.56%  "

$ws.Range("E43").Value = "  -2.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9999"
$ws.Range("E44").Value = "  +0.02%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.727"
$ws.Range("E45").Value = "  +2.21%  "

$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000123"
$ws.Range("E46").Value = "  +2.16%  "

$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "1.991.84"
$ws.Range("E47").Value = "  -1.07%  "

$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5186"
$ws.Range("E48").Value = "  +0.09%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "64.38"
$ws.Range("E49").Value = "  -0.74%  "

$ws.Range("E50").Value = "  +1.22%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.009"
$ws.Range("E51").Value = "  -0.01%  "
